# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-18 17:13:22
#
# The "Recorded By" column (G) stores a comma-separated list of recorders
# for each session row. A handful of distinct values had their first
# entry rotated to the end of the list. Apply the same literal text
# replacements to every matching cell in column G on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colG = $ws.Range("G1:G" + $ws.UsedRange.Rows.Count)

$colG.Replace("system, backup@backdoor.com, System", "backup@backdoor.com, System, system")
$colG.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System")
$colG.Replace("admin@admin.com, dnasr281@gmail.com", "dnasr281@gmail.com, admin@admin.com")
